$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data for rows 2..17 (A index, B name, C from_bus, D to_bus, E in_service)
# Two new "line7"/"line8" rows are inserted (sorted) right after "line6", pushing the
# previously-existing "extr1".."extr8" rows down by two, and all rows get refreshed
# from_bus / to_bus / in_service values per the latest contingency run.

$data = @(
    @(0,  "line1", 7,  9,  $true),
    @(1,  "line2", 9,  8,  $true),
    @(2,  "line3", 8,  10, $false),
    @(3,  "line4", 8,  11, $true),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $true),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $true),
    @(11, "extr4", 7,  8,  $true),
    @(12, "extr5", 9,  11, $false),
    @(13, "extr6", 7,  11, $true),
    @(14, "extr7", 5,  7,  $false),
    @(15, "extr8", 8,  5,  $false)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}

# Rows 16 and 17 are brand new rows; copy the column-A formatting (bold, border,
# centered) used by the other index cells (style index 1) onto them.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
